# Insert a new weekly data row for "Terminal La Palmera de La Serena - Cilantro"
# at row 59, pushing the existing rows 59..137 down to 60..138.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 59 (shifts rows 59-137 down to 60-138)
$ws.Rows.Item(59).Insert()

# Populate the new row with the latest week's values
$ws.Cells.Item(59, 1).Value = 8
$ws.Cells.Item(59, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44665
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = 100112040
$ws.Cells.Item(59, 7).Value = "Cilantro"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 3000
$ws.Cells.Item(59, 11).Value = 2000
$ws.Cells.Item(59, 12).Value = 2500
$ws.Cells.Item(59, 13).Value = 2250
$ws.Cells.Item(59, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(59, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(59, 16).Value = 1500
$ws.Cells.Item(59, 17).Value = 1.5
$ws.Cells.Item(59, 18).Value = "Hortaliza"
